$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Strip the "Mazda_RX8_Coupe_1/" prefix from the block-name labels in A2:A28
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $s = $cell.Text
    $cell.Value = ($s -replace "Mazda_RX8_Coupe_1/", "")
}

# 2. Collapse the redundant duplicate style used by B15/B18 back onto the
#    shared default style (copy formats from a cell that already uses it).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 3. Tweak the tab-bar/scrollbar split ratio shown in the window chrome.
$excel.ActiveWindow.TabRatio = 0.983

# 4. Move the active cell to A1 while keeping the rest of the column selected.
$r1 = $ws.Range("A2:A28")
$r2 = $ws.Range("A1")
$u = $excel.Union($r1, $r2)
$u.Select() | Out-Null

# 5. Column A got a touch wider after the labels shrank - nudge its width.
$ws.Columns.Item(1).ColumnWidth = 36.8333333
